$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the first data row (old row 2, "H 72"), shifting all rows up by one
$ws.Rows.Item(2).Delete()

# Apply value corrections for cells where the missing-data mask differs from the shifted source
$ws.Range("D2").Value = -14.8
$ws.Range("E2").Value = -6.9
$ws.Range("D3").Value = ""
$ws.Range("E4").Value = ""
$ws.Range("E6").Value = -6.4
$ws.Range("E8").Value = ""
$ws.Range("D13").Value = -14.7
$ws.Range("D14").Value = ""
$ws.Range("D15").Value = -15.4
$ws.Range("D16").Value = ""
$ws.Range("E16").Value = -5.3
$ws.Range("E18").Value = ""
$ws.Range("D20").Value = -15.3
$ws.Range("E20").Value = -6.9
$ws.Range("D21").Value = ""
$ws.Range("D22").Value = -15.2
$ws.Range("E22").Value = ""
$ws.Range("D23").Value = ""
$ws.Range("D24").Value = -14
$ws.Range("D26").Value = ""
$ws.Range("D28").Value = -13.9
$ws.Range("E28").Value = -7
$ws.Range("D29").Value = ""
$ws.Range("E32").Value = ""
$ws.Range("E37").Value = -7.1
$ws.Range("D39").Value = -14.7
$ws.Range("E39").Value = ""
$ws.Range("D40").Value = ""
$ws.Range("E40").Value = -7.9
$ws.Range("E42").Value = ""
$ws.Range("D52").Value = -13.8
$ws.Range("D53").Value = ""
$ws.Range("D56").Value = -14.7
$ws.Range("E56").Value = -5.7
$ws.Range("D57").Value = ""
$ws.Range("D58").Value = -13
$ws.Range("E58").Value = ""
$ws.Range("D59").Value = ""
$ws.Range("E62").Value = -10.7

# Re-affirm already-blank cells so they are stored consistently as empty (no residual shared-string "")
$ws.Range("D25").Value = ""
$ws.Range("E30").Value = ""
